$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Enquiry" columns/values describing the Unit Linked Mortality Charges param
$ws.Range("Z2").Value = "Unit Linked Mortality Charges based in Age"
$ws.Range("AA2").Value = "Unit Linked Mortality Charges for Male and Female Lives"
$ws.Range("Z3").Value = "Rate of Unit Linked Mortality Charge"

# Shared formatting (small font, vertical-centered wrapped text) across the whole block
$all = $ws.Range("Z2:AA3")
$all.Font.Size = 7.5
$all.VerticalAlignment = -4108
$all.WrapText = $true

# Thin border around Z2:AA2 and Z3 (AA3 is left without a border)
$ws.Range("Z2:AA2").Borders.LineStyle = 1
$ws.Range("Z3").Borders.LineStyle = 1

# Column widths for the new columns (closest values achievable given this
# engine's pixel-quantized column-width storage)
$ws.Columns.Item(26).ColumnWidth = 48.3
$ws.Columns.Item(27).ColumnWidth = 39.0

# Selection / view, matching the saved workbook state
$null = $ws.Range("Z2:AA3").Select()

Write-Output "done"
